# diary.xlsx update — "First aggregate demo skeleton in plac"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: extend the 28 loka time range and log the ImGui work ---
# New "Huomiot koodista" entry for the 28 loka row, describing the
# evening ImGui integration work. (entered before the Kello edit below)
$ws.Range("F17").Value = "ImGui integrointi iltapuhteena, tästä tulee hyvä"
$ws.Range("F17").WrapText = $true

# Kello: 8.45-10.45, 11.45-13.15  ->  8.45-10.45, 11.45-13.15, 19.45-20.45
$ws.Range("B17").Value = "8.45-10.45, 11.45-13.15, 19.45-20.45"

# Hours logged for 28 loka grows from 3.5 to 4.5 (the extra 19.45-20.45 hour)
$ws.Range("G17").Value = 4.5

# --- Row 18 (new): 29 loka ---
$ws.Range("A18").Value = "29 loka"

$ws.Range("C18").Value = "Kytketyn kappaleen demoa"
$ws.Range("C18").WrapText = $true

$ws.Range("B18").Value = "19.15-21.45"
$ws.Range("B18").NumberFormat = "h:mm"

$ws.Range("D18").Value = "Jospa nyt olisi se particle.cpp integrointi metodi kunnossa :D . Vielä pitää korjata firework, mutta eka kytketty kappale demo pohja valmiina."
$ws.Range("D18").WrapText = $true

$ws.Range("G18").Value = 2.5

$ws.Rows.Item(18).RowHeight = 43.5

# --- View state: scrolled down one row further, selection moved to D19 ---
$ws.Range("D19").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
